# Car Test Transaction Specification.docx - apply commit "inserts and specs changed"
#
# The edit touches the title (re-themed + "s" added -> "Specifications"),
# inserts a new bold "Transaction:" summary paragraph (with a blank bold-flagged
# paragraph above it), indents + substantially rewrites the "Description:"
# paragraph, and appends a trailing "Enter a test no. ..." paragraph.
#
# Because the new/rewritten runs need precise run-level formatting (theme
# fonts, run-level <w:sz>, a <w:proofErr> spell-check bookmark pair, an
# indent-only paragraph, and an entirely empty paragraph) that line up
# exactly with the target markup, the whole body is re-expressed as literal
# WordprocessingML and dropped in through Range.InsertXML - the supported
# COM path for "replace this range's contents with exact OOXML" (confirmed
# by the host itself: Range.WordOpenXML is read-only and directs XML edits
# through InsertXML).

$d = $word.ActiveDocument

$newBodyXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Car Test Transaction Specification</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>s</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Transaction</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> Enter the date to check if any cars are required to take a test due to cars age or date since last test. If they do require test, create a Letter for the Owner</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Description:</w:t></w:r><w:r><w:t xml:space="preserve"> The transactions purpose is to take in a</w:t></w:r><w:r><w:t xml:space="preserve"> date</w:t></w:r><w:r><w:t xml:space="preserve"> through a front end and </w:t></w:r><w:r><w:t xml:space="preserve">return the owners of cars that are older than 5 years or cars that are older than 5 years and haven&#8217;t had a test in a year, </w:t></w:r><w:r><w:t xml:space="preserve">by joining the Car table with the Owner and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CarRecord</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> table</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>If those conditions are true, a letter</w:t></w:r><w:r><w:t xml:space="preserve"> will be created for the owners</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> The transaction will also look to see if a letter with a recent date exists for that owner, to stop a new letter being created every day after they due to be tested. </w:t></w:r><w:r><w:t>There will be an exception catc</w:t></w:r><w:r><w:t>h for the possibility that no cars are due to be tested that day</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Enter a test no. and check if there are any high criticality failures</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Content.InsertXML($newBodyXml)
